$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.021786099299788475
$ws.Range("C2").Value = 0.00978164840489626
$ws.Range("D2").Value = 0.007128716912120581
$ws.Range("E2").Value = 0.006709726061671972
$ws.Range("F2").Value = 0.00026360771153122187
$ws.Range("J2").Value = 0.12743115425109863
$ws.Range("K2").Value = 1.4351342916488647
